# 20240720 Update - Scrape duration
# Overall update to the code so it is usable again and improved the README.md
#
# Rework the sample sheet: replace the old wide "one row per location-type"
# layout with a simple Origin/Destination matrix header (A1:C1) plus a
# single column (A2:A15) listing all the stations / landmarks used by the
# scraper, including the newly added locations (浅草寺, 上野駅, 豊洲駅,
# 東京タワー, 東京スカイツリー). Stale entries (Origin, Beemars,
# Harry's Sandwich, 十条駅, the Shinjuku address, ...) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the stale D1:O1 header values (Origin/address/Beemars/etc.), but
#     keep the two special cell styles already sitting on G1 (fontId 3 /
#     細明體) and O1 (fontId 2 / Arial) in place. ---
$ws.Range("D1:O1").ClearContents()

# --- Row 1, column A: new Origin/Destination matrix header ---
$ws.Cells.Item(1, 1).Value = "Origin - Destination"

# --- Column A: the list of stations / landmarks (rows 2-15) ---
$values = @(
  "中目黒駅",
  "渋谷駅",
  "代官山駅",
  "恵比寿駅",
  "麻布十番駅",
  "六本木駅",
  "新宿駅",
  "飯田橋駅",
  "日本橋駅",
  "銀座駅",
  "大手町駅",
  "浅草寺",
  "上野駅",
  "豊洲駅"
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $values[$i]
}

# A7 and A15 reuse the same 細明體 style that G1 already carries, so copy
# G1's format onto them (rather than setting Font.Name directly, which
# would mint a brand-new style entry).
$ws.Cells.Item(1, 7).Copy() | Out-Null
$ws.Cells.Item(7, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 1, columns B & C: filled in last ---
$ws.Cells.Item(1, 2).Value = "東京タワー"
$ws.Cells.Item(1, 3).Value = "東京スカイツリー"

# --- Selection, matching the saved workbook view ---
$ws.Range("E13").Select() | Out-Null
